$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6198
$ws.Range('K3').Value = 6395
$ws.Range('K4').Value = 1335
$ws.Range('K5').Value = 456
$ws.Range('K6').Value = 7044
$ws.Range('K7').Value = 21428
$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('K2').Value = 3
$ws.Range('K6').Value = 22
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K3').Value = 73
$ws.Range('K6').Value = 115
$ws.Range('K7').Value = 273
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 430
$ws.Range('K4').Value = 77
$ws.Range('K6').Value = 475
$ws.Range('K7').Value = 1404
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 162
$ws.Range('K7').Value = 467
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 242
$ws.Range('K3').Value = 336
$ws.Range('K6').Value = 285
$ws.Range('K7').Value = 934
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K3').Value = 124
$ws.Range('K7').Value = 353
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 243
$ws.Range('K6').Value = 211
$ws.Range('K7').Value = 725
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 125
$ws.Range('K6').Value = 179
$ws.Range('K7').Value = 504
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 145
$ws.Range('K7').Value = 352
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K7').Value = 629
$ws.Range('K8').Value = 1404
$ws.Range('K9').Value = 92
$ws.Range('K10').Value = 123
$ws.Range('K11').Value = 400
$ws.Range('K13').Value = 30
$ws.Range('K14').Value = 106
$ws.Range('K18').Value = 141
$ws.Range('K20').Value = 512
$ws.Range('K24').Value = 66
$ws.Range('K29').Value = 1166
$ws.Range('K33').Value = 934
$ws.Range('K34').Value = 123
$ws.Range('K37').Value = 725
$ws.Range('K38').Value = 22
$ws.Range('K39').Value = 27
$ws.Range('K42').Value = 794
$ws.Range('K48').Value = 269
$ws.Range('K49').Value = 117
$ws.Range('K51').Value = 276
$ws.Range('K52').Value = 565
$ws.Range('K53').Value = 273
$ws.Range('K54').Value = 418
$ws.Range('K55').Value = 233
$ws.Range('K62').Value = 8
$ws.Range('K63').Value = 58
$ws.Range('K64').Value = 134
$ws.Range('K65').Value = 504
$ws.Range('K67').Value = 841
$ws.Range('K71').Value = 64
$ws.Range('K73').Value = 190
$ws.Range('K76').Value = 290
$ws.Range('K77').Value = 147
$ws.Range('K79').Value = 540
$ws.Range('K83').Value = 467
$ws.Range('K85').Value = 999
$ws.Range('K87').Value = 39
$ws.Range('K89').Value = 316
$ws.Range('K90').Value = 198
$ws.Range('K91').Value = 246
$ws.Range('K95').Value = 353
$ws.Range('K96').Value = 224
$ws.Range('K97').Value = 170
$ws.Range('K99').Value = 352
$ws.Range('K101').Value = 21428
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K6').Value = 239
$ws.Range('K7').Value = 841
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 23
$ws.Range('K7').Value = 117
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 100
$ws.Range('K6').Value = 229
$ws.Range('K7').Value = 418
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 332
$ws.Range('K3').Value = 418
$ws.Range('K6').Value = 331
$ws.Range('K7').Value = 1166
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K6').Value = 128
$ws.Range('K7').Value = 269
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 66
$ws.Range('K5').Value = 3
$ws.Range('K7').Value = 290
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K3').Value = 22
$ws.Range('K7').Value = 106
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 240
$ws.Range('K7').Value = 794
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K3').Value = 10
$ws.Range('K6').Value = 30
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K2').Value = 38
$ws.Range('K7').Value = 123
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 80
$ws.Range('K7').Value = 233
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K2').Value = 27
$ws.Range('K3').Value = 14
$ws.Range('K7').Value = 66
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K6').Value = 97
$ws.Range('K7').Value = 224
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 119
$ws.Range('K7').Value = 246
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 175
$ws.Range('K6').Value = 136
$ws.Range('K7').Value = 540
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K5').Value = 3
$ws.Range('K7').Value = 134
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 175
$ws.Range('K3').Value = 167
$ws.Range('K7').Value = 512
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 141
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 205
$ws.Range('K3').Value = 208
$ws.Range('K5').Value = 24
$ws.Range('K6').Value = 170
$ws.Range('K7').Value = 629
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K6').Value = 36
$ws.Range('K7').Value = 123
$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('K5').Value = 17
$ws.Range('K6').Value = 27
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 141
$ws.Range('K5').Value = 4
$ws.Range('K6').Value = 129
$ws.Range('K7').Value = 400
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K2').Value = 28
$ws.Range('K7').Value = 92
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 63
$ws.Range('K7').Value = 190
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K2').Value = 34
$ws.Range('K7').Value = 170
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 96
$ws.Range('K7').Value = 316
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 74
$ws.Range('K3').Value = 57
$ws.Range('K7').Value = 198
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K6').Value = 90
$ws.Range('K7').Value = 276
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 326
$ws.Range('K5').Value = 30
$ws.Range('K6').Value = 246
$ws.Range('K7').Value = 999
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K3').Value = 22
$ws.Range('K7').Value = 64
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K6').Value = 18
$ws.Range('K7').Value = 147
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 164
$ws.Range('K7').Value = 565
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('K2').Value = 5
$ws.Range('K3').Value = 11
$ws.Range('K7').Value = 39
$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range('K6').Value = 6
$ws.Range('K7').Value = 8
